$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 49, shifting existing rows 49:77 down to 50:78.
# Use xlFormatFromRightOrBelow (-4161) so the new row inherits the date-style
# formatting (style index 2) from the row below it, just like all the other
# data rows.
$ws.Rows.Item(49).Insert(-4161)

# Fill the new row 49 with the new data record (same categorical values as the
# surrounding rows, new date/price values per the diff)
$ws.Cells.Item(49, 1).Value = 9
$ws.Cells.Item(49, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(49, 3).Value = "Metropolitana"
$ws.Cells.Item(49, 4).Value = 44813
$ws.Cells.Item(49, 5).Value = 13
$ws.Cells.Item(49, 6).Value = 100112029
$ws.Cells.Item(49, 7).Value = "Orégano"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 16
$ws.Cells.Item(49, 11).Value = 18000
$ws.Cells.Item(49, 12).Value = 18000
$ws.Cells.Item(49, 13).Value = 18000
$ws.Cells.Item(49, 14).Value = "$/docena de atados"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 6000
$ws.Cells.Item(49, 17).Value = 3
$ws.Cells.Item(49, 18).Value = "Hortaliza"
